# Corrected minor typos on W9S3.
#
# Slide 77 ("If time allows, let us define our Tokenizer v1.0"), in the
# Content Placeholder shape, has a paragraph reading:
#   "Will classify the lexemes, one at a time, and will assemble create
#    a Token struct for each lexeme, which will contain:"
# "will assemble create" is a typo for "will create". Fix it by
# replacing that substring in place (this is what causes PowerPoint to
# split the paragraph's single run into three runs at the edited span).

$oldFragment = "will assemble create "
$newFragment = "will create "

$p = $ppt.ActivePresentation

$fixed = $false

for ($si = 1; $si -le $p.Slides.Count -and -not $fixed; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count -and -not $fixed; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if (-not $shp.HasTextFrame) { continue }

        $tr = $shp.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count

        for ($pi = 1; $pi -le $paraCount -and -not $fixed; $pi++) {
            $para = $tr.Paragraphs($pi)
            $fullText = $para.Text
            $startIdx = $fullText.IndexOf($oldFragment)

            if ($startIdx -ge 0) {
                $target = $para.Characters($startIdx + 1, $oldFragment.Length)
                $target.Text = $newFragment
                $fixed = $true
            }
        }
    }
}
